$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting for the new column F from existing columns so the
# existing (header / plain-border) styles are reused rather than new
# style entries being minted.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("D2:D9").Copy()
$ws.Range("F2:F9").PasteSpecial(-4122)

# New header + values in column F
$ws.Range("F1").Value = "Jira"
$ws.Range("F2").Value = "Qa-0124"
$ws.Range("F6").Value = "Qa-0125"

# Match the recorded selection left behind in the worksheet view
$ws.Range("M17").Select() | Out-Null
